$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("target_data")
$ws.Range("A1").Value = "test"
